$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear K3:M5 and K12:M14 (these rows no longer have K/L/M data)
$ws.Range("K3:M5").ClearContents()
$ws.Range("K12:M14").ClearContents()

# Update K6:M11 with new values
$ws.Range("K6").Value = 112
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 4

$ws.Range("K7").Value = 180
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 4

$ws.Range("K8").Value = 140
$ws.Range("L8").Value = 5
$ws.Range("M8").Value = 4

$ws.Range("K9").Value = 308
$ws.Range("L9").Value = 11
$ws.Range("M9").Value = 4

$ws.Range("K10").Value = 532
$ws.Range("L10").Value = 19
$ws.Range("M10").Value = 4

$ws.Range("K11").Value = 812
$ws.Range("L11").Value = 29
$ws.Range("M11").Value = 4

# Add new values to B41:D41
$ws.Range("B41").Value = 812
$ws.Range("C41").Value = 29
$ws.Range("D41").Value = 4

# Update the selected range/active cell
$ws.Range("M11").Select()
